$d = $word.ActiveDocument

# --- Edit 1: "ARROW FUNCTION" heading gains a trailing (bold) space before
#     " adalah bentuk lain ..." ---
$r1 = $d.Content
$r1.Find.ClearFormatting()
$r1.Find.Execute("ARROW FUNCTION", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "ARROW FUNCTION ", 2)

# --- Edit 2: the "Contoh 5" object-literal line collapses its split runs
#     back down to the same visible code (no textual change), e.g.
#     "let jumlahHuruf = mahasiswa.map ( nama =>" + " ({" + " " -> one run,
#     and "nama.length" + "}" + ");" -> one run. ---
$r2 = $d.Content
$r2.Find.ClearFormatting()
$r2.Find.Execute("let jumlahHuruf = mahasiswa.map ( nama => ({ nama: nama, jmlHuruf: nama.length});", `
                  $true, $false, $false, $false, $false, `
                  $true, 1, $false, `
                  "let jumlahHuruf = mahasiswa.map ( nama => ({ nama: nama, jmlHuruf: nama.length});", 2)

Write-Output "done"
